$wb = $excel.ActiveWorkbook

# 1. Update the "In Translation" status text everywhere it is used
#    (Overview!E/F, zh-cn!C, de-de!C all share this text).
foreach ($ws in $wb.Worksheets) {
    $ws.Cells.Replace("In Translation", "Handed back: in sync with en-US")
}

# 2. Fill in the "Latest Target File" / "Latest Handback File" / "Latest
#    Handback DateTime" columns now that the handback has happened.

$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$md1 = "df76b354-7d8e-4cd0-8cdf-b75543d3360d.md"
$md2 = "fdb77385-ecad-4252-91e7-1dad5dc2591d.md"
$md1url = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/69641fae25987c8bef373b223407f6c0f4a2ce19/e2e/df76b354-7d8e-4cd0-8cdf-b75543d3360d.md"
$md2url = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/69641fae25987c8bef373b223407f6c0f4a2ce19/e2e/fdb77385-ecad-4252-91e7-1dad5dc2591d.md"

# zh-cn sheet, row 2 (df76b354...) and row 3 (fdb77385...)
$zhcn.Range("I2").Value = $md1
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), $md1url, "", "", $md1) | Out-Null
$zhcn.Range("J2").Value = "df76b354-7d8e-4cd0-8cdf-b75543d3360d.4526fd6d5cc5aeccf86b7cee0d39c31ae9fa31ef.zh-cn.xlf"

$zhcn.Range("I3").Value = $md2
$zhcn.Hyperlinks.Add($zhcn.Range("I3"), $md2url, "", "", $md2) | Out-Null
$zhcn.Range("J3").Value = "fdb77385-ecad-4252-91e7-1dad5dc2591d.8c5bd0a6d7205f3a02ffc540137660b3a6b5e0dd.zh-cn.xlf"

# de-de sheet, row 2 (df76b354...) and row 3 (fdb77385...)
$dede.Range("I2").Value = $md1
$dede.Hyperlinks.Add($dede.Range("I2"), $md1url, "", "", $md1) | Out-Null
$dede.Range("J2").Value = "df76b354-7d8e-4cd0-8cdf-b75543d3360d.4526fd6d5cc5aeccf86b7cee0d39c31ae9fa31ef.de-de.xlf"
$dede.Range("K2").Value = "2016-08-29 08:26:36"

$dede.Range("I3").Value = $md2
$dede.Hyperlinks.Add($dede.Range("I3"), $md2url, "", "", $md2) | Out-Null
$dede.Range("J3").Value = "fdb77385-ecad-4252-91e7-1dad5dc2591d.8c5bd0a6d7205f3a02ffc540137660b3a6b5e0dd.de-de.xlf"
$dede.Range("K3").Value = "2016-08-29 08:26:36"

# zh-cn "Latest Handback DateTime" column (K) now populated too.
$zhcn.Range("K2").Value = "2016-08-29 08:26:29"
$zhcn.Range("K3").Value = "2016-08-29 08:26:29"

# 3. Widen the columns that now hold longer file-name / status text.
$wb.Worksheets.Item("Overview").Range("E1").ColumnWidth = 29.9777047293527
$wb.Worksheets.Item("Overview").Range("F1").ColumnWidth = 29.9777047293527

foreach ($ws in @($zhcn, $dede)) {
    $ws.Range("C1").ColumnWidth = 29.9777047293527
    $ws.Range("I1").ColumnWidth = 40
    $ws.Range("J1").ColumnWidth = 40
}
